# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型"
# sheets to reflect the latest scrape (both rows now read 383).

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 383
    $ws.Range("F3").Value = 383
}
